$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 10 run times (Part 1 / Part 2) - fills in the previously blank cells
$ws.Range("B14").Value = 0.014606699987780299
$ws.Range("C14").Value = 0.072334899974521194

# Recalculate so dependent formulas (E14, B31, C31, E31) update
$excel.Calculate()

# Move the active selection to E14, matching the saved view state
$ws.Range("E14").Select()
